$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing data row (row 2) values into the new row 3
$ws.Range("A2:H2").Copy()
$ws.Range("A3").PasteSpecial(-4163)

# Column I gets a new value for the new row
$ws.Range("I3").Value = "SIN_ASIGNAR, adios"
